# Applies the "update database and change read_price algorithm" edit:
# the 5-year rolling income-statement window shifts left by one column
# (D<-E, E<-F, F<-G, G<-H) and a new year of figures lands in column H,
# for both header rows (period / publish-date labels) and the data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "دوره مالی" (financial period) column headers, D8:H8 ---
$periodHeaders = @("12 ماهه منتهی به 1397/12", "12 ماهه منتهی به 1398/12", "12 ماهه منتهی به 1399/12", "12 ماهه منتهی به 1400/12", "12 ماهه منتهی به 1401/12")
$cols = @("D","E","F","G","H")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periodHeaders[$i]
}

# --- Row 9: "تاریخ انتشار" (publish date) column headers, D9:H9 ---
$dateHeaders = @("1399-02-22 (12)", "1400-03-02 (14)", "1401-02-26 (9)", "1402-02-29 (8)", "1402-02-29")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $dateHeaders[$i]
}

# --- Data rows 11-37, columns D:H ---
# Each entry: row number, then the 5 values for D,E,F,G,H (use "-" for the dash placeholder)
$rowsData = @(
    @(11, "-", "-", 1794029, 3123661, 6662826),
    @(12, "-", "-", 92184, 107297, 186247),
    @(13, "-", "-", 1894273, 3240025, 6624368),
    @(14, "-", "-", -1871825, -3198990, -6597688),
    @(15, 806540, 1182073, 1908661, 3271993, 6875753),
    @(16, "-", "-", -1446944, -2625935, -5909215),
    @(17, -48182, -62464, -93357, -184457, -341706),
    @(18, 9189, 10240, 20458, 25482, 21460),
    @(19, -8336, -26627, -29017, -37305, -55766),
    @(20, 205007, 292857, 359801, 449778, 590526),
    @(21, "-", "-", 55182, 108091, 177755),
    @(22, 1742, 108, 1375, 1897, 3824),
    @(23, "-", "-", "-", 4964, -6140),
    @(24, -12488, -159, 0, -12371, -92835),
    @(25, 254812, 333302, 417768, 552359, 673130),
    @(26, -121, 0, -317, 0, -646),
    @(27, 0, -216, 0, -740, -1698),
    @(28, 254691, 333086, 417451, 551619, 670786),
    @(29, 0, 0, 0, 0, 0),
    @(30, 254691, 333086, 417451, 551619, 670786),
    @(31, 0, 0, 0, 0, 0),
    @(32, 0, 0, 0, 0, 0),
    @(33, 0, 0, 0, 0, 0),
    @(34, 0, 0, 0, 0, 0),
    @(35, 1273, 1665, 2087, 788, 958),
    @(36, 200000, 200000, 200000, 700000, 700000),
    @(37, 364, 476, 596, 788, 958)
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $entry[$i + 1]
    }
}
